$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(5, 6).Value = 133
$ws.Cells.Item(6, 6).Value = 10428
$ws.Cells.Item(6, 7).Value = 120
$ws.Cells.Item(8, 6).Value = 3636
$ws.Cells.Item(10, 6).Value = 2478
$ws.Cells.Item(12, 6).Value = 2910
$ws.Cells.Item(14, 6).Value = 517
$ws.Cells.Item(15, 6).Value = 2234
$ws.Cells.Item(19, 6).Value = 409
$ws.Cells.Item(24, 6).Value = 270
$ws.Cells.Item(25, 6).Value = 624
$ws.Cells.Item(28, 6).Value = 1276
$ws.Cells.Item(32, 6).Value = 3995
$ws.Cells.Item(33, 6).Value = 3481
$ws.Cells.Item(36, 6).Value = 1070
$ws.Cells.Item(37, 6).Value = 421
$ws.Cells.Item(40, 6).Value = 124

# --- Sheet "本地生活" (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(5, 6).Value = 2154

# --- Sheet "全部类型" (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(9, 6).Value = 133
$ws.Cells.Item(10, 6).Value = 10428
$ws.Cells.Item(10, 7).Value = 120
$ws.Cells.Item(12, 6).Value = 3636
$ws.Cells.Item(14, 6).Value = 2478
$ws.Cells.Item(16, 6).Value = 2910
$ws.Cells.Item(17, 6).Value = 517
$ws.Cells.Item(18, 6).Value = 2234
$ws.Cells.Item(22, 6).Value = 409
$ws.Cells.Item(25, 6).Value = 270
$ws.Cells.Item(26, 6).Value = 624
$ws.Cells.Item(29, 6).Value = 1276
$ws.Cells.Item(33, 6).Value = 3995
$ws.Cells.Item(34, 6).Value = 3481
$ws.Cells.Item(36, 6).Value = 1070
$ws.Cells.Item(38, 6).Value = 421
$ws.Cells.Item(44, 6).Value = 124
